# Add a new "AutoAtkDis" property/column to the Skill sheet.
# This inserts a new column I (between the existing AtkDis and NeedTar
# columns), shifting NeedTar -> J and DefaultHitTime -> K, then fills in
# the new header and a default value of 1 for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (9th column), pushing the old
# NeedTar / DefaultHitTime columns one slot to the right.
$ws.Columns("I:I").Insert()

# New column header.
$ws.Range("I1").Value = "AutoAtkDis"

# Default value for the new property on every existing data row (2-9).
$lastRow = $ws.Cells(1, 1).End(-4121).Row
if ($lastRow -lt 2) { $lastRow = 9 }
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
}

# Match the selection left behind by the edit in the source workbook.
$ws.Range("I2:I9").Select() | Out-Null
